# Update cryptocurrency price/volume data in the active worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($addr, $text) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

Set-TextCell "D2" "42.037.11"
Set-TextCell "E2" "  -0.44%  "
Set-TextCell "D3" "2.188.38"
Set-TextCell "E3" "  -2.57%  "
Set-TextCell "E4" "  -0.17%  "
Set-TextCell "D5" "239.15"
Set-TextCell "E5" "  -1.76%  "
Set-TextCell "D6" "0.607"
Set-TextCell "E6" "  -3.15%  "
Set-TextCell "D7" "73.30"
Set-TextCell "E7" "  -1.47%  "
Set-TextCell "D9" "0.588"
Set-TextCell "E9" "  -2.43%  "
Set-TextCell "D10" "40.43"
Set-TextCell "E10" "  -4.48%  "
Set-TextCell "E11" "  -4.48%  "
Set-TextCell "D12" "54.72"
Set-TextCell "E12" "  -3.42%  "
Set-TextCell "D13" "6.80"
Set-TextCell "E13" "  -2.36%  "
Set-TextCell "E14" "  -2.71%  "
Set-TextCell "D15" "2.512.38"
Set-TextCell "E15" "  -2.73%  "
Set-TextCell "D16" "14.52"
Set-TextCell "E16" "  +0.86%  "
Set-TextCell "D17" "2.180.44"
Set-TextCell "E17" "  -2.43%  "
Set-TextCell "D18" "0.788"
Set-TextCell "E18" "  -6.41%  "
Set-TextCell "D19" "41.844.43"
Set-TextCell "E19" "  -0.73%  "
Set-TextCell "D20" "0.0000103"
Set-TextCell "E20" "  -1.85%  "
Set-TextCell "D21" "70.38"
Set-TextCell "E21" "  -3.36%  "
Set-TextCell "E22" "  -6.53%  "
Set-TextCell "D23" "10.13"
Set-TextCell "E23" "  -9.42%  "
Set-TextCell "D24" "227.40"
Set-TextCell "E24" "  -1.37%  "
Set-TextCell "D25" "2.07"
Set-TextCell "E25" "  +1.27%  "
Set-TextCell "E26" "  -0.06%  "
Set-TextCell "D27" "10.80"
Set-TextCell "E27" "  -5.88%  "
Set-TextCell "E28" "  -9.33%  "
Set-TextCell "E29" "  -3.37%  "
Set-TextCell "E30" "  -1.01%  "
Set-TextCell "D31" "172.00"
Set-TextCell "E31" "  +2.46%  "
Set-TextCell "D32" "20.01"
Set-TextCell "E32" "  -3.04%  "
Set-TextCell "D33" "33.10"
Set-TextCell "E33" "  +10.36%  "
Set-TextCell "D34" "0.0784"
Set-TextCell "E34" "  -2.69%  "
Set-TextCell "D35" "5.32"
Set-TextCell "E35" "  -6.18%  "
Set-TextCell "E36" "  -3.22%  "
Set-TextCell "D37" "4.39"
Set-TextCell "E37" "  +1.28%  "
Set-TextCell "E38" "  -5.81%  "
Set-TextCell "D39" "0.0314"
Set-TextCell "E39" "  +2.65%  "
Set-TextCell "D40" "12.34"
Set-TextCell "E40" "  -6.44%  "
Set-TextCell "E41" "  -1.82%  "
Set-TextCell "D42" "5.41"
Set-TextCell "E42" "  -5.47%  "
Set-TextCell "D43" "59.81"
Set-TextCell "E43" "  -8.07%  "
Set-TextCell "E44" "  -3.72%  "
Set-TextCell "D46" "0.0970"
Set-TextCell "E46" "  -3.75%  "
Set-TextCell "D47" "98.29"
Set-TextCell "E47" "  -5.95%  "
Set-TextCell "E48" "  -4.51%  "
Set-TextCell "D49" "1.13"
Set-TextCell "E49" "  -4.30%  "
Set-TextCell "D50" "2.23"
Set-TextCell "E50" "  -5.03%  "
Set-TextCell "D51" "0.418"
Set-TextCell "E51" "  +12.92%  "
